$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$logRows = @(
    ,@("2025-12-27 13:22:09", "Admin", "Login", "login_success", "Role: admin")
    ,@("2025-12-27 13:22:09", "Admin", "dashboard", "access_granted", "Opened dashboard page")
    ,@("2025-12-27 13:22:19", "Admin", "dashboard", "access_granted", "Opened dashboard page")
    ,@("2025-12-27 13:22:35", "Admin", "quotation", "access_granted", "Opened quotation page")
    ,@("2025-12-27 13:23:08", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:23:11", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:23:11", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:23:14", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:23:20", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:23:21", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:23:32", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:23:36", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:23:39", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:23:39", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:23:41", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:23:41", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:23:43", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:24:19", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:24:20", "Admin", "receipt", "access_granted", "Opened receipt page")
    ,@("2025-12-27 13:24:21", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:24:23", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:24:28", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:24:32", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:24:35", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:24:35", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:24:37", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:25:16", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:25:18", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:25:32", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:25:36", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:28:27", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:28:45", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:33:27", "Admin", "Login", "login_success", "Role: admin")
    ,@("2025-12-27 13:33:27", "Admin", "dashboard", "access_granted", "Opened dashboard page")
    ,@("2025-12-27 13:33:30", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:33:34", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:33:35", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:33:36", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:33:42", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:33:47", "Admin", "invoice", "access_granted", "Opened invoice page")
    ,@("2025-12-27 13:33:53", "Admin", "invoice", "access_granted", "Opened invoice page")
)

$startRow = 90
for ($i = 0; $i -lt $logRows.Count; $i++) {
    $r = $startRow + $i
    $data = $logRows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
}

"Done: added rows 90-130"
